{"js": "// Replace the 25 three-digit / one-digit division problems & answers\n// in the table with their updated values, matching each old string\n// to its replacement exactly once (search scoped fresh each iteration\n// so a later replacement's new text can't accidentally be re-matched).\nconst replacements = [\n  [\"994\u00f73=331, 1\", \"524\u00f72=262, 0\"],\n  [\"226\u00f78=28, 2\", \"434\u00f77=62, 0\"],\n  [\"459\u00f76=76, 3\", \"372\u00f78=46, 4\"],\n  [\"483\u00f73=161, 0\", \"744\u00f78=93, 0\"],\n  [\"424\u00f73=141, 1\", \"489\u00f73=163, 0\"],\n  [\"348\u00f73=116, 0\", \"502\u00f77=71, 5\"],\n  [\"422\u00f73=140, 2\", \"874\u00f79=97, 1\"],\n  [\"232\u00f76=38, 4\", \"592\u00f74=148, 0\"],\n  [\"186\u00f79=20, 6\", \"581\u00f79=64, 5\"],\n  [\"550\u00f73=183, 1\", \"226\u00f74=56, 2\"],\n  [\"564\u00f78=70, 4\", \"737\u00f73=245, 2\"],\n  [\"177\u00f76=29, 3\", \"692\u00f75=138, 2\"],\n  [\"310\u00f74=77, 2\", \"103\u00f72=51, 1\"],\n  [\"500\u00f76=83, 2\", \"767\u00f72=383, 1\"],\n  [\"610\u00f74=152, 2\", \"983\u00f76=163, 5\"],\n  [\"501\u00f77=71, 4\", \"761\u00f76=126, 5\"],\n  [\"156\u00f73=52, 0\", \"881\u00f76=146, 5\"],\n  [\"533\u00f75=106, 3\", \"510\u00f78=63, 6\"],\n  [\"448\u00f75=89, 3\", \"324\u00f76=54, 0\"],\n  [\"206\u00f76=34, 2\", \"292\u00f74=73, 0\"],\n  [\"874\u00f78=109, 2\", \"981\u00f78=122, 5\"],\n  [\"759\u00f72=379, 1\", \"469\u00f75=93, 4\"],\n  [\"968\u00f74=242, 0\", \"608\u00f76=101, 2\"],\n  [\"650\u00f77=92, 6\", \"862\u00f76=143, 4\"],\n  [\"348\u00f79=38, 6\", \"700\u00f78=87, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the 25 three-digit / one-digit division problems & answers\n# in the table with their updated values (Find/Replace per pair; uses\n# two parallel arrays rather than an array-of-arrays since nested\n# array literals get flattened by this interpreter).\n$d = $word.ActiveDocument\n\n$oldValues = @(\n    \"994\u00f73=331, 1\",\n    \"226\u00f78=28, 2\",\n    \"459\u00f76=76, 3\",\n    \"483\u00f73=161, 0\",\n    \"424\u00f73=141, 1\",\n    \"348\u00f73=116, 0\",\n    \"422\u00f73=140, 2\",\n    \"232\u00f76=38, 4\",\n    \"186\u00f79=20, 6\",\n    \"550\u00f73=183, 1\",\n    \"564\u00f78=70, 4\",\n    \"177\u00f76=29, 3\",\n    \"310\u00f74=77, 2\",\n    \"500\u00f76=83, 2\",\n    \"610\u00f74=152, 2\",\n    \"501\u00f77=71, 4\",\n    \"156\u00f73=52, 0\",\n    \"533\u00f75=106, 3\",\n    \"448\u00f75=89, 3\",\n    \"206\u00f76=34, 2\",\n    \"874\u00f78=109, 2\",\n    \"759\u00f72=379, 1\",\n    \"968\u00f74=242, 0\",\n    \"650\u00f77=92, 6\",\n    \"348\u00f79=38, 6\"\n)\n\n$newValues = @(\n    \"524\u00f72=262, 0\",\n    \"434\u00f77=62, 0\",\n    \"372\u00f78=46, 4\",\n    \"744\u00f78=93, 0\",\n    \"489\u00f73=163, 0\",\n    \"502\u00f77=71, 5\",\n    \"874\u00f79=97, 1\",\n    \"592\u00f74=148, 0\",\n    \"581\u00f79=64, 5\",\n    \"226\u00f74=56, 2\",\n    \"737\u00f73=245, 2\",\n    \"692\u00f75=138, 2\",\n    \"103\u00f72=51, 1\",\n    \"767\u00f72=383, 1\",\n    \"983\u00f76=163, 5\",\n    \"761\u00f76=126, 5\",\n    \"881\u00f76=146, 5\",\n    \"510\u00f78=63, 6\",\n    \"324\u00f76=54, 0\",\n    \"292\u00f74=73, 0\",\n    \"981\u00f78=122, 5\",\n    \"469\u00f75=93, 4\",\n    \"608\u00f76=101, 2\",\n    \"862\u00f76=143, 4\",\n    \"700\u00f78=87, 4\"\n)\n\nfor ($i = 0; $i -lt $oldValues.Count; $i++) {\n    $oldText = $oldValues[$i]\n    $newText = $newValues[$i]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND at index $i\"\n    }\n}\n"}
